$wb = $excel.ActiveWorkbook

# Budget sheet: the budget year in A2 moves from 2021 to 2022
$budget = $wb.Worksheets.Item("Budget")
$budget.Range("A2").Value = 2022

# Budget becomes the active/selected tab (was Criteria before), with A2 selected
$budget.Select()
$budget.Range("A2").Select()

$wb.Save()
